$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-02-14 Wednesday" "2024-02-15 Thursday"

Replace-Text "97÷6=" "58÷5="
Replace-Text "98÷2=" "61÷5="
Replace-Text "31÷7=" "16÷6="
Replace-Text "67÷2=" "37÷6="
Replace-Text "32÷7=" "96÷4="

Replace-Text "65÷9=" "50÷6="
Replace-Text "93÷8=" "62÷6="
Replace-Text "68÷9=" "51÷6="
Replace-Text "13÷3=" "70÷6="
Replace-Text "58÷2=" "62÷9="

Replace-Text "58÷3=" "46÷4="
Replace-Text "83÷9=" "97÷2="
Replace-Text "69÷9=" "24÷3="
Replace-Text "75÷2=" "78÷5="
Replace-Text "47÷3=" "16÷5="

Replace-Text "14÷4=" "46÷7="
Replace-Text "90÷6=" "57÷5="
Replace-Text "45÷3=" "25÷5="
Replace-Text "60÷9=" "14÷5="
Replace-Text "87÷4=" "46÷4="

Replace-Text "82÷6=" "99÷2="
Replace-Text "91÷3=" "58÷6="
Replace-Text "46÷9=" "88÷7="
Replace-Text "96÷7=" "59÷2="
Replace-Text "14÷8=" "27÷9="
